$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hungary NB I")

# --- Update existing rows 150 and 151 with new match data, and add new odds columns ---
# Row 150
$ws.Range("A150").Value = 148
$ws.Range("B150").Value = 6818351
$ws.Range("C150").Value = "Hungary NB I"
$ws.Range("D150").Value = "Hungary NB I"
$ws.Range("E150").Value = 45380.66666666666
$ws.Range("F150").Value = "Puskas Academy"
$ws.Range("G150").Value = "MOL Fehervar FC"
$ws.Range("H150").Value = 0
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = "D"
$ws.Range("K150").Value = 2.45
$ws.Range("L150").Value = 3.3
$ws.Range("M150").Value = 2.45
$ws.Range("N150").Value = 1.65
$ws.Range("O150").Value = 3.75
$ws.Range("P150").Value = 4.75
$ws.Range("Q150").Value = -0.75
$ws.Range("R150").Value = 1.875
$ws.Range("S150").Value = 1.975
$ws.Range("T150").Value = 2.5
$ws.Range("U150").Value = 1.875
$ws.Range("V150").Value = 1.975
$ws.Range("W150").Value = -1
$ws.Range("X150").Value = 2.75
$ws.Range("Y150").Value = -1
$ws.Range("Z150").Value = -1
$ws.Range("AA150").Value = 0.9750000000000001
$ws.Range("AB150").Value = -1
$ws.Range("AC150").Value = 0.9750000000000001

# Row 151
$ws.Range("A151").Value = 149
$ws.Range("B151").Value = 6818347
$ws.Range("C151").Value = "Hungary NB I"
$ws.Range("D151").Value = "Hungary NB I"
$ws.Range("E151").Value = 45381.4375
$ws.Range("F151").Value = "Kisvarda FC"
$ws.Range("G151").Value = "Debreceni VSC"
$ws.Range("H151").Value = 1
$ws.Range("I151").Value = 3
$ws.Range("J151").Value = "A"
$ws.Range("K151").Value = 2.75
$ws.Range("L151").Value = 3.25
$ws.Range("M151").Value = 2.25
$ws.Range("N151").Value = 2.875
$ws.Range("O151").Value = 3.3
$ws.Range("P151").Value = 2.15
$ws.Range("Q151").Value = 0.25
$ws.Range("R151").Value = 1.825
$ws.Range("S151").Value = 2.025
$ws.Range("T151").Value = 2.25
$ws.Range("U151").Value = 1.875
$ws.Range("V151").Value = 1.975
$ws.Range("W151").Value = -1
$ws.Range("X151").Value = -1
$ws.Range("Y151").Value = 1.15
$ws.Range("Z151").Value = -1
$ws.Range("AA151").Value = 1.025
$ws.Range("AB151").Value = 0.875
$ws.Range("AC151").Value = -1

# --- Copy style (border/font/alignment for col A, date format for col E) from template row 149 to new rows 152-158 ---
$ws.Range("A149").Copy() | Out-Null
$ws.Range("A152:A158").PasteSpecial(-4122) | Out-Null
$ws.Range("E149").Copy() | Out-Null
$ws.Range("E152:E158").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Add new rows 152-158 (new matches) ---
# Row 152
$ws.Range("A152").Value = 150
$ws.Range("B152").Value = 6818346
$ws.Range("C152").Value = "Hungary NB I"
$ws.Range("D152").Value = "Hungary NB I"
$ws.Range("E152").Value = 45381.54166666666
$ws.Range("F152").Value = "MTK Budapest"
$ws.Range("G152").Value = "Kecskemeti TE"
$ws.Range("H152").Value = 2
$ws.Range("I152").Value = 2
$ws.Range("J152").Value = "D"
$ws.Range("K152").Value = 2.3
$ws.Range("L152").Value = 3.25
$ws.Range("M152").Value = 2.7
$ws.Range("N152").Value = 2.1
$ws.Range("O152").Value = 3.3
$ws.Range("P152").Value = 3
$ws.Range("Q152").Value = -0.25
$ws.Range("R152").Value = 1.95
$ws.Range("S152").Value = 1.9
$ws.Range("T152").Value = 2.5
$ws.Range("U152").Value = 2
$ws.Range("V152").Value = 1.85
$ws.Range("W152").Value = -1
$ws.Range("X152").Value = 2.3
$ws.Range("Y152").Value = -1
$ws.Range("Z152").Value = -0.5
$ws.Range("AA152").Value = 0.45
$ws.Range("AB152").Value = 1
$ws.Range("AC152").Value = -1

# Row 153
$ws.Range("A153").Value = 151
$ws.Range("B153").Value = 6818348
$ws.Range("C153").Value = "Hungary NB I"
$ws.Range("D153").Value = "Hungary NB I"
$ws.Range("E153").Value = 45381.64583333334
$ws.Range("F153").Value = "Zalaegerszegi TE"
$ws.Range("G153").Value = "Diosgyori VTK"
$ws.Range("H153").Value = 5
$ws.Range("I153").Value = 1
$ws.Range("J153").Value = "H"
$ws.Range("K153").Value = 2.45
$ws.Range("L153").Value = 3.3
$ws.Range("M153").Value = 2.45
$ws.Range("N153").Value = 2.5
$ws.Range("O153").Value = 3.3
$ws.Range("P153").Value = 2.4
$ws.Range("Q153").Value = 0
$ws.Range("R153").Value = 1.975
$ws.Range("S153").Value = 1.875
$ws.Range("T153").Value = 2.75
$ws.Range("U153").Value = 2.025
$ws.Range("V153").Value = 1.825
$ws.Range("W153").Value = 1.5
$ws.Range("X153").Value = -1
$ws.Range("Y153").Value = -1
$ws.Range("Z153").Value = 0.9750000000000001
$ws.Range("AA153").Value = -1
$ws.Range("AB153").Value = 1.025
$ws.Range("AC153").Value = -1

# Row 154
$ws.Range("A154").Value = 152
$ws.Range("B154").Value = 6818353
$ws.Range("C154").Value = "Hungary NB I"
$ws.Range("D154").Value = "Hungary NB I"
$ws.Range("E154").Value = 45387.625
$ws.Range("F154").Value = "Ujpest"
$ws.Range("G154").Value = "Zalaegerszegi TE"
$ws.Range("K154").Value = 1.909
$ws.Range("L154").Value = 3.5
$ws.Range("M154").Value = 3.25
$ws.Range("N154").Value = 2
$ws.Range("O154").Value = 3.4
$ws.Range("P154").Value = 3.3
$ws.Range("Q154").Value = -0.25
$ws.Range("R154").Value = 1.85
$ws.Range("S154").Value = 2
$ws.Range("T154").Value = 2.5
$ws.Range("U154").Value = 1.9
$ws.Range("V154").Value = 1.95
$ws.Range("W154").Value = 0
$ws.Range("X154").Value = 0
$ws.Range("Y154").Value = 0
$ws.Range("Z154").Value = 0
$ws.Range("AA154").Value = 0

# Row 155
$ws.Range("A155").Value = 153
$ws.Range("B155").Value = 6818356
$ws.Range("C155").Value = "Hungary NB I"
$ws.Range("D155").Value = "Hungary NB I"
$ws.Range("E155").Value = 45388.38541666666
$ws.Range("F155").Value = "MOL Fehervar FC"
$ws.Range("G155").Value = "Mezokovesd Zsory"
$ws.Range("K155").Value = 1.571
$ws.Range("L155").Value = 3.6
$ws.Range("M155").Value = 5
$ws.Range("N155").Value = 1.727
$ws.Range("O155").Value = 3.4
$ws.Range("P155").Value = 4.5
$ws.Range("Q155").Value = -0.75
$ws.Range("R155").Value = 2.025
$ws.Range("S155").Value = 1.825
$ws.Range("T155").Value = 2.5
$ws.Range("U155").Value = 1.85
$ws.Range("V155").Value = 2
$ws.Range("W155").Value = 0
$ws.Range("X155").Value = 0
$ws.Range("Y155").Value = 0
$ws.Range("Z155").Value = 0
$ws.Range("AA155").Value = 0

# Row 156
$ws.Range("A156").Value = 154
$ws.Range("B156").Value = 6818354
$ws.Range("C156").Value = "Hungary NB I"
$ws.Range("D156").Value = "Hungary NB I"
$ws.Range("E156").Value = 45388.47916666666
$ws.Range("F156").Value = "Diosgyori VTK"
$ws.Range("G156").Value = "Kisvarda FC"
$ws.Range("K156").Value = 1.727
$ws.Range("L156").Value = 3.4
$ws.Range("M156").Value = 4.2
$ws.Range("N156").Value = 1.8
$ws.Range("O156").Value = 3.3
$ws.Range("P156").Value = 4.2
$ws.Range("Q156").Value = -0.5
$ws.Range("R156").Value = 1.875
$ws.Range("S156").Value = 1.975
$ws.Range("T156").Value = 2.5
$ws.Range("U156").Value = 1.825
$ws.Range("V156").Value = 2.025
$ws.Range("W156").Value = 0
$ws.Range("X156").Value = 0
$ws.Range("Y156").Value = 0
$ws.Range("Z156").Value = 0
$ws.Range("AA156").Value = 0

# Row 157
$ws.Range("A157").Value = 155
$ws.Range("B157").Value = 6818355
$ws.Range("C157").Value = "Hungary NB I"
$ws.Range("D157").Value = "Hungary NB I"
$ws.Range("E157").Value = 45388.60416666666
$ws.Range("F157").Value = "Debreceni VSC"
$ws.Range("G157").Value = "MTK Budapest"
$ws.Range("K157").Value = 2
$ws.Range("L157").Value = 3.4
$ws.Range("M157").Value = 3.1
$ws.Range("N157").Value = 1.8
$ws.Range("O157").Value = 3.5
$ws.Range("P157").Value = 4
$ws.Range("Q157").Value = -0.5
$ws.Range("R157").Value = 1.875
$ws.Range("S157").Value = 1.975
$ws.Range("T157").Value = 2.5
$ws.Range("U157").Value = 1.825
$ws.Range("V157").Value = 2.025
$ws.Range("W157").Value = 0
$ws.Range("X157").Value = 0
$ws.Range("Y157").Value = 0
$ws.Range("Z157").Value = 0
$ws.Range("AA157").Value = 0

# Row 158
$ws.Range("A158").Value = 156
$ws.Range("B158").Value = 6818352
$ws.Range("C158").Value = "Hungary NB I"
$ws.Range("D158").Value = "Hungary NB I"
$ws.Range("E158").Value = 45389.40625
$ws.Range("F158").Value = "Ferencvarosi TC"
$ws.Range("G158").Value = "Paksi"
$ws.Range("K158").Value = 1.444
$ws.Range("L158").Value = 4
$ws.Range("M158").Value = 5.75
$ws.Range("N158").Value = 1.45
$ws.Range("O158").Value = 4
$ws.Range("P158").Value = 6.5
$ws.Range("Q158").Value = -1
$ws.Range("R158").Value = 1.8
$ws.Range("S158").Value = 2.05
$ws.Range("T158").Value = 3
$ws.Range("U158").Value = 1.9
$ws.Range("V158").Value = 1.95
$ws.Range("W158").Value = 0
$ws.Range("X158").Value = 0
$ws.Range("Y158").Value = 0
$ws.Range("Z158").Value = 0
$ws.Range("AA158").Value = 0
